# Improve doco, rename mass to body mass, add min & max trait values from ranges
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the "mass" column definitions to "body mass" (values only - row
#    numbers are unaffected since nothing is inserted yet).
$ws.Range("A14").Value = "body mass"
$ws.Range("A15").Value = "body mass - units"
$ws.Range("A16").Value = "body mass - method"
$ws.Range("A17").Value = "body mass - comments"
$ws.Range("A18").Value = "body mass - metadata comment"
$ws.Range("A19").Value = "original body mass"
$ws.Range("A20").Value = "original body mass - units"

# 2) Insert the two new "body mass" range rows right after "body mass - units" (row 15).
$ws.Rows.Item(16).Insert()
$ws.Rows.Item(16).Insert()
$ws.Range("A16").Value = "body mass - minimum"
$ws.Range("B16").Value = "Minimum body mass (in standardised units) if a range was specified in the data source"
$ws.Range("A17").Value = "body mass - maximum"
$ws.Range("B17").Value = "Maximum body mass (in standardised units) if a range was specified in the data source"

# 3) Insert the two new "metabolic rate" range rows right after "metabolic rate - units" (now row 24).
$ws.Rows.Item(25).Insert()
$ws.Rows.Item(25).Insert()
$ws.Range("A25").Value = "metabolic rate - minimum"
$ws.Range("A26").Value = "metabolic rate - maximum"
$ws.Range("B25").Value = "Minimum metabolic rate (in standardised units) if a range was specified in the data source"
$ws.Range("B26").Value = "Maximum metabolic rate (in standardised units) if a range was specified in the data source"

# 4) Insert the two new "brain size" range rows right after "brain size - units" (now row 37).
$ws.Rows.Item(38).Insert()
$ws.Rows.Item(38).Insert()
$ws.Range("A38").Value = "brain size - minimum"
$ws.Range("A39").Value = "brain size - maximum"
$ws.Range("B38").Value = "Minimum brain size (in standardised units) if a range was specified in the data source"
$ws.Range("B39").Value = "Maximum brain size (in standardised units) if a range was specified in the data source"

# 5) Fix up the view so it matches the author's saved cursor/selection position.
$ws.Range("A36:XFD39").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 11
$win.ScrollColumn = 1
